$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 476 entirely; all subsequent rows (477-493) shift up by one
# (becoming 476-492), matching the diff which removes the
# "珈琲で一番大事な事！誰と飲むか？" post row and renumbers the rest.
$ws.Rows.Item(476).Delete()
